# P6-3-4-BLECapsense.docx edit script
# Applies the "first crack at 3-4" content rewrite.

$d = $word.ActiveDocument
$rq = [char]0x2019   # right single quotation mark used throughout the doc
$hellip = [char]0x2026

# ---------------------------------------------------------------------
# 1) Paragraph 2 ("Welcome back ...") - trim the run-on sentence after
#    "... for our robotic arm."
# ---------------------------------------------------------------------
$old1 = "for our robotic arm. So, we" + $rq + "ll start a new project, and add a CapSense capacitive-sensing interface to control the robotic arm that is connected to the PSoC 6 kit configured as a peripheral device, via BLE."
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, "for our robotic arm.", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Paragraph 3 ("To get started, ...") -> "Add the capsense component"
# ---------------------------------------------------------------------
$old2 = "To get started, let" + $rq + "s create a new PSoC Creator project, we" + $rq + "ll call it BLE CapSense Remote Control"
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, "Add the capsense component", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Paragraph 4 (yellow) -> "Change the name to capsense"
# ---------------------------------------------------------------------
$old3 = "[Create a new project, add and configure the CapSense Component, show the PDL APIs]"
$d.Content.Find.Execute($old3, $false, $false, $false, $false, $false, $true, 1, $false, "Change the name to capsense", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) Paragraph 5 (yellow) -> "Add linear slider and two buttons"
# ---------------------------------------------------------------------
$old4 = "[Add and configure the BLE Component, show the PDL APIs]"
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, "Add linear slider and two buttons", 2) | Out-Null

# ---------------------------------------------------------------------
# 5) Paragraph 6 (yellow) -> "Make the buttons mutual cap"
# ---------------------------------------------------------------------
$old5 = "[Explain that for this project, we" + $rq + "ll again dedicate the CM0+ for the BLE functionality and leave the CM4 to do the CapSense functions and what" + $rq + "s to come]"
$d.Content.Find.Execute($old5, $false, $false, $false, $false, $false, $true, 1, $false, "Make the buttons mutual cap", 2) | Out-Null

# ---------------------------------------------------------------------
# 6) Paragraph 7 (yellow) -> "Go to the advanced tab -> widget details... change the Button1_tx to be button0_tx"
# ---------------------------------------------------------------------
$old6 = "[Add and describe the firmware across the two cores]"
$new6 = "Go to the advanced tab -> widget details" + $hellip + " change the Button1_tx to be button0_tx"
$d.Content.Find.Execute($old6, $false, $false, $false, $false, $false, $true, 1, $false, $new6, 2) | Out-Null

# ---------------------------------------------------------------------
# 7) Paragraph 8 (" [Build and run]", yellow) -> pin assignment paragraph
# ---------------------------------------------------------------------
$old7 = " [Build and run]"
$new7 = "Assign the pins the linear slider to pin P8[3] -> P8[7], The Button0 RX to P8[1] and Button 1 to P8[2] " + $hellip + " then the Button Tx to P1[0] " + $hellip + " then you assign the capacitors to their default location" + $hellip + " see they are labled in green when you do the pulldown menu"
$d.Content.Find.Execute($old7, $false, $false, $false, $false, $false, $true, 1, $false, $new7, 2) | Out-Null

# ---------------------------------------------------------------------
# 8) Paragraph 9 (Demo..., yellow) -> code-reuse paragraph
# ---------------------------------------------------------------------
$old8 = "[Demo and show how to connect the two PSoC 6 BLE kits and control the robotic arm with the CapSense interface]"
$new8 = "Now" + $hellip + " because I am into code reuse Ill copy the capsenseTask.h and .c from the MainController project" + $hellip + " lets see I need to set the includes to be project.h freertos.h task.h and "
$d.Content.Find.Execute($old8, $false, $false, $false, $false, $false, $true, 1, $false, $new8, 2) | Out-Null

# ---------------------------------------------------------------------
# 9) Move the _GoBack bookmark from the "Now we have..." paragraph to the
#    end of the paragraph we just rewrote (paragraph 9), and insert a new
#    blank paragraph after it.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
if ($bm.Exists) { $bm.Delete() }

$p9 = $d.Paragraphs.Item(9)
$endOfP9 = $p9.Range.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($endOfP9, $endOfP9)) | Out-Null

$p9.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------
# 10) Paragraph "Now we have..." sentence tweak
# ---------------------------------------------------------------------
$old9 = "lets add in some sensors. For the next few videos, we" + $rq + "ll be implementing the motion sensor and temperature sensor on the E-ink Display shield board to the BLE remote controller!"
$new9 = "lets add in some sensors. For the next few videos, we" + $rq + "ll be implementing the motion sensor and the E-ink Display to the BLE remote controller!"
$d.Content.Find.Execute($old9, $false, $false, $false, $false, $false, $true, 1, $false, $new9, 2) | Out-Null

# ---------------------------------------------------------------------
# 11) Strip the now-unused yellow highlighting throughout the document.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.HighlightColorIndex -ne 0) {
        $p.Range.Style = "Normal (Web)"
        $p.Range.Font.NameAscii = "Arial"
        $p.Range.Font.Name = "Arial"
        $p.Range.Font.NameBi = "Arial"
        $p.Range.HighlightColorIndex = 0
    }
}
